# Auto-generated: update Ixion_Profits market-data values per commit diff.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1239.4445
$ws.Range("I62").Value = 1244.375
$ws.Range("K62").Value = 1244.375
$ws.Range("M62").Value = -620.375
$ws.Range("H65").Value = 1239.4445
$ws.Range("I65").Value = 1244.375
$ws.Range("K65").Value = 6221.875
$ws.Range("M65").Value = -3101.875
$ws.Range("H123").Value = 28664.455
$ws.Range("J123").Value = 28664.455
$ws.Range("L123").Value = 28664.455
$ws.Range("N123").Value = -38464.455
$ws.Range("H137").Value = 1710.5217
$ws.Range("I137").Value = 1442.1
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 4326.299999999999
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -1776.299999999999
$ws.Range("N137").Value = -15600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4138.9116
$ws.Range("I32").Value = 3024.3584
$ws.Range("J32").Value = 8077
$ws.Range("K32").Value = 3024.3584
$ws.Range("L32").Value = 8077
$ws.Range("M32").Value = -2737.3584
$ws.Range("N32").Value = -8651
$ws.Range("H61").Value = 4919.0605
$ws.Range("I61").Value = 5446.5713
$ws.Range("J61").Value = 1965
$ws.Range("K61").Value = 5446.5713
$ws.Range("L61").Value = 1965
$ws.Range("M61").Value = -5234.5713
$ws.Range("N61").Value = -2389
$ws.Range("H74").Value = 1630.0807
$ws.Range("I74").Value = 1588.8596
$ws.Range("K74").Value = 1588.8596
$ws.Range("M74").Value = -714.8596
$ws.Range("H77").Value = 1630.0807
$ws.Range("I77").Value = 1588.8596
$ws.Range("K77").Value = 7944.298
$ws.Range("M77").Value = -3576.298
$ws.Range("H132").Value = 3693.5312
$ws.Range("I132").Value = 1807.0625
$ws.Range("J132").Value = 5580
$ws.Range("K132").Value = 5421.1875
$ws.Range("L132").Value = 16740
$ws.Range("M132").Value = -2891.1875
$ws.Range("N132").Value = -21800
$ws.Range("H136").Value = 4919.0605
$ws.Range("I136").Value = 5446.5713
$ws.Range("J136").Value = 1965
$ws.Range("K136").Value = 16339.7139
$ws.Range("L136").Value = 5895
$ws.Range("M136").Value = -13789.7139
$ws.Range("N136").Value = -10995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4923.3335
$ws.Range("I134").Value = 5861.04
$ws.Range("J134").Value = 2792.182
$ws.Range("K134").Value = 17583.12
$ws.Range("L134").Value = 8376.545999999998
$ws.Range("M134").Value = -15048.12
$ws.Range("N134").Value = -13446.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4505.2046
$ws.Range("I31").Value = 1208.5172
$ws.Range("J31").Value = 10878.8
$ws.Range("K31").Value = 1208.5172
$ws.Range("L31").Value = 10878.8
$ws.Range("M31").Value = -913.5172
$ws.Range("N31").Value = -11468.8
$ws.Range("H34").Value = 4505.2046
$ws.Range("I34").Value = 1208.5172
$ws.Range("J34").Value = 10878.8
$ws.Range("K34").Value = 1208.5172
$ws.Range("L34").Value = 10878.8
$ws.Range("M34").Value = -1006.5172
$ws.Range("N34").Value = -11282.8
$ws.Range("H58").Value = 1594.4517
$ws.Range("I58").Value = 1400
$ws.Range("J58").Value = 1717.2632
$ws.Range("K58").Value = 1400
$ws.Range("L58").Value = 1717.2632
$ws.Range("M58").Value = -1197
$ws.Range("N58").Value = -2123.2632
$ws.Range("H99").Value = 11368181
$ws.Range("I99").Value = 2448.6667
$ws.Range("K99").Value = 2448.6667
$ws.Range("M99").Value = -950.6667000000002
$ws.Range("H126").Value = 11368181
$ws.Range("I126").Value = 2448.6667
$ws.Range("K126").Value = 7346.000100000001
$ws.Range("M126").Value = -4876.000100000001
$ws.Range("H132").Value = 2270.775
$ws.Range("I132").Value = 2046.3438
$ws.Range("J132").Value = 3168.5
$ws.Range("K132").Value = 6139.0314
$ws.Range("L132").Value = 9505.5
$ws.Range("M132").Value = -3609.0314
$ws.Range("N132").Value = -14565.5
$ws.Range("H134").Value = 1572.68
$ws.Range("I134").Value = 1471.35
$ws.Range("K134").Value = 4414.049999999999
$ws.Range("M134").Value = -1879.049999999999
$ws.Range("H136").Value = 1594.4517
$ws.Range("I136").Value = 1400
$ws.Range("J136").Value = 1717.2632
$ws.Range("K136").Value = 4200
$ws.Range("L136").Value = 5151.7896
$ws.Range("M136").Value = -1650
$ws.Range("N136").Value = -10251.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 5882537
$ws.Range("I12").Value = 11111282
$ws.Range("J12").Value = 198.75
$ws.Range("K12").Value = 33333846
$ws.Range("L12").Value = 596.25
$ws.Range("M12").Value = -33333673
$ws.Range("N12").Value = -942.25
$ws.Range("H34").Value = 2583.1667
$ws.Range("J34").Value = 3574.75
$ws.Range("L34").Value = 10724.25
$ws.Range("N34").Value = -10892.25
$ws.Range("H60").Value = 451.14285
$ws.Range("I60").Value = 431.6
$ws.Range("J60").Value = 500
$ws.Range("K60").Value = 1294.8
$ws.Range("L60").Value = 1500
$ws.Range("M60").Value = -1043.8
$ws.Range("N60").Value = -2002
$ws.Range("H62").Value = 6075
$ws.Range("J62").Value = 6075
$ws.Range("L62").Value = 18225
$ws.Range("N62").Value = -19597
$ws.Range("H65").Value = 6075
$ws.Range("J65").Value = 6075
$ws.Range("L65").Value = 54675
$ws.Range("N65").Value = -61539
$ws.Range("H131").Value = 1695886.5
$ws.Range("I131").Value = 7143650
$ws.Range("J131").Value = 1026.8667
$ws.Range("K131").Value = 21430950
$ws.Range("L131").Value = 3080.6001
$ws.Range("M131").Value = -21425910
$ws.Range("N131").Value = -13160.6001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5160.76
$ws.Range("I132").Value = 5265.7
$ws.Range("J132").Value = 5090.8
$ws.Range("K132").Value = 15797.1
$ws.Range("L132").Value = 15272.4
$ws.Range("M132").Value = -13267.1
$ws.Range("N132").Value = -20332.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3204.4546
$ws.Range("I61").Value = 2874.8333
$ws.Range("K61").Value = 2874.8333
$ws.Range("M61").Value = -2672.8333
$ws.Range("H113").Value = 3204.4546
$ws.Range("I113").Value = 2874.8333
$ws.Range("K113").Value = 2874.8333
$ws.Range("M113").Value = -704.8332999999998
$ws.Range("H132").Value = 9634474
$ws.Range("I132").Value = 12386304
$ws.Range("K132").Value = 37158912
$ws.Range("M132").Value = -37156382
$ws.Range("H136").Value = 4868.8726
$ws.Range("I136").Value = 4241.2
$ws.Range("J136").Value = 6699.5835
$ws.Range("K136").Value = 12723.6
$ws.Range("L136").Value = 20098.7505
$ws.Range("M136").Value = -10173.6
$ws.Range("N136").Value = -25198.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 36333.332
$ws.Range("I123").Value = 30000
$ws.Range("J123").Value = 49000
$ws.Range("K123").Value = 30000
$ws.Range("L123").Value = 49000
$ws.Range("M123").Value = -25100
$ws.Range("N123").Value = -58800
$ws.Range("H131").Value = 44000
$ws.Range("J131").Value = 44000
$ws.Range("L131").Value = 44000
$ws.Range("N131").Value = -54080
$ws.Range("H132").Value = 1526.7576
$ws.Range("I132").Value = 966.4286
$ws.Range("J132").Value = 2507.3333
$ws.Range("K132").Value = 2899.2858
$ws.Range("L132").Value = 7521.999899999999
$ws.Range("M132").Value = -369.2857999999997
$ws.Range("N132").Value = -12581.9999
$ws.Range("H136").Value = 2316.5789
$ws.Range("I136").Value = 2354.1428
$ws.Range("K136").Value = 7062.428400000001
$ws.Range("M136").Value = -4512.428400000001
